# Consolidate the "CORDSET" category label into the existing "CORD-SET"
# label (rows 38-47, 53, 98-100 on Sheet1), retarget row 89's stray
# "PLAZO" category to "SKIRT-TOP", and rename row 102's category from
# "CORD-SET" to the newly introduced "CORD-SUIT" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 38-47: Category column (B) "CORDSET" -> "CORD-SET"
foreach ($r in 38..47) {
    $ws.Cells.Item($r, 2).Value = "CORD-SET"
}

# Row 53: Category column (B) "CORDSET" -> "CORD-SET"
$ws.Cells.Item(53, 2).Value = "CORD-SET"

# Rows 98-100: Category column (B) "CORDSET" -> "CORD-SET"
foreach ($r in 98..100) {
    $ws.Cells.Item($r, 2).Value = "CORD-SET"
}

# Row 89: Category column (B) "PLAZO" -> "SKIRT-TOP"
$ws.Cells.Item(89, 2).Value = "SKIRT-TOP"

# Row 102: Category column (B) "CORD-SET" -> "CORD-SUIT"
$ws.Cells.Item(102, 2).Value = "CORD-SUIT"
